$d = $word.ActiveDocument

# --- Header 1: "1. Will Holt - Monday, [Date], [Time]" ---
$p = $d.Paragraphs.Item(2)
$rng = $p.Range
$rng.Find.Execute("1. Will Holt - Monday, [Date], [Time]", $true, $false, $false, $false, $false, $true, 0, $false, "1. Will Holt – Monday 3/3 11:05am", 2)

# --- Header 2: "2. Vaishnavi Paniki - Monday, [Date], [Time]" (only trailing run) ---
$p = $d.Paragraphs.Item(11)
$rng = $p.Range
$rng.Find.Execute(" - Monday, [Date], [Time]", $true, $false, $false, $false, $false, $true, 0, $false, " – Monday 2/2 11:30am", 2)

# --- Header 3: "3. Will Holt - Monday, [Date], [Time]" ---
$p = $d.Paragraphs.Item(18)
$rng = $p.Range
$rng.Find.Execute("3. Will Holt - Monday, [Date], [Time]", $true, $false, $false, $false, $false, $true, 0, $false, "3. Will Holt – Tuesday 3/4 10:19am", 2)

# --- Header 4: "4. Sebastian Segura - Monday, [Date], [Time]" ---
$p = $d.Paragraphs.Item(26)
$rng = $p.Range
$rng.Find.Execute("4. Sebastian Segura - Monday, [Date], [Time]", $true, $false, $false, $false, $false, $true, 0, $false, "4. Sebastian Segura – Wednesday 3/5 4:51pm", 2)

# --- Header 5: "5. Will Holt - Monday, [Date], [Time]" ---
$p = $d.Paragraphs.Item(32)
$rng = $p.Range
$rng.Find.Execute("5. Will Holt - Monday, [Date], [Time]", $true, $false, $false, $false, $false, $true, 0, $false, "5. Will Holt – Wednesday 3/5 6:03pm", 2)

# --- Header 6: "6. Will Holt - [Date], [Time]" ---
$p = $d.Paragraphs.Item(37)
$rng = $p.Range
$rng.Find.Execute("6. Will Holt - [Date], [Time]", $true, $false, $false, $false, $false, $true, 0, $false, "6. Will Holt – Friday 3/7 2:29pm", 2)

# --- Header 7: "7. Vaishnavi Paniki - [Date], [Time]" (only trailing run) ---
$p = $d.Paragraphs.Item(49)
$rng = $p.Range
$rng.Find.Execute(" - [Date], [Time]", $true, $false, $false, $false, $false, $true, 0, $false, " – Saturday 3/8 11:28pm", 2)

# --- Remove the trailing "Let me know..." paragraph and the blank paragraph after it ---
$count = $d.Paragraphs.Count
$pRef = $d.Paragraphs.Item($count - 1)
$pLast = $d.Paragraphs.Item($count)
$delRange = $d.Range($pRef.Range.Start, $pLast.Range.End)
$delRange.Delete()

Write-Output "done"
